# Refresh crypto price/volume snapshot values (and swap the two rows whose
# market-cap rank changed) to match the latest scrape, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: some "Price" values look like plain decimals (e.g. "0.490", "7.10")
# which Excel would otherwise auto-convert to numbers on assignment, silently
# dropping the trailing zero / formatting. Force those specific cells to Text
# format first so the literal string is preserved, matching the sheet convention
# of storing these as text. Values that are not parseable as a single plain
# number (e.g. thousand-grouped "63.630.20", or percentages like "  -3.54%  ")
# are assigned directly since Excel keeps those as text already.

# Row 2
$ws.Range("D2").Value2 = "63.630.20"
$ws.Range("E2").Value2 = "  -3.54%  "

# Row 3
$ws.Range("D3").Value2 = "3.481.04"
$ws.Range("E3").Value2 = "  -2.98%  "

# Row 4
$ws.Range("E4").Value2 = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.16"
$ws.Range("E5").Value2 = "  -3.01%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.97"
$ws.Range("E6").Value2 = "  -5.65%  "

# Row 7
$ws.Range("D7").Value2 = "3.481.65"
$ws.Range("E7").Value2 = "  -2.94%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.490"
$ws.Range("E9").Value2 = "  -1.87%  "

# Row 10
$ws.Range("E10").Value2 = "  -1.73%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.10"
$ws.Range("E11").Value2 = "  -1.57%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.382"
$ws.Range("E12").Value2 = "  -2.63%  "

# Row 13
$ws.Range("D13").Value2 = "4.075.12"
$ws.Range("E13").Value2 = "  -3.00%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.60"
$ws.Range("E14").Value2 = "  -2.05%  "

# Row 15
$ws.Range("E15").Value2 = "  +1.05%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000178"
$ws.Range("E16").Value2 = "  -4.49%  "

# Row 17
$ws.Range("D17").Value2 = "3.492.70"

# Row 18
$ws.Range("D18").Value2 = "63.770.41"
$ws.Range("E18").Value2 = "  -3.43%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.07"
$ws.Range("E19").Value2 = "  +0.21%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.34"
$ws.Range("E20").Value2 = "  -1.84%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.66"
$ws.Range("E21").Value2 = "  -3.20%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "387.63"
$ws.Range("E22").Value2 = "  -2.39%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.576"
$ws.Range("E23").Value2 = "  -1.97%  "

# Row 24
$ws.Range("D24").Value2 = "3.625.17"
$ws.Range("E24").Value2 = "  -2.93%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.03"
$ws.Range("E25").Value2 = "  -2.66%  "

# Row 26
$ws.Range("E26").Value2 = "  +0.08%  "

# Row 27
$ws.Range("E27").Value2 = "  -6.65%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.59"
$ws.Range("E28").Value2 = "  -3.05%  "

# Row 29
$ws.Range("B29").Value2 = "Binance-PegBSC-USD"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value2 = "  -0.14%  "

# Row 30
$ws.Range("B30").Value2 = "RenderToken"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.41"
$ws.Range("E30").Value2 = "  -8.57%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.24"
$ws.Range("E31").Value2 = "  -3.94%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.19"
$ws.Range("E32").Value2 = "  -5.20%  "

# Row 33
$ws.Range("D33").Value2 = "3.486.96"
$ws.Range("E33").Value2 = "  -3.00%  "

# Row 34
$ws.Range("E34").Value2 = "  +0.01%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.73"
$ws.Range("E35").Value2 = "  -3.24%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.143"
$ws.Range("E36").Value2 = "  -4.36%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.27"
$ws.Range("E37").Value2 = "  -2.58%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.56"
$ws.Range("E38").Value2 = "  -3.11%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.87"
$ws.Range("E39").Value2 = "  -2.21%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "167.53"
$ws.Range("E40").Value2 = "  -0.76%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0801"
$ws.Range("E41").Value2 = "  -4.46%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.76"
$ws.Range("E42").Value2 = "  +3.06%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.809"
$ws.Range("E43").Value2 = "  -3.70%  "

# Row 44
$ws.Range("E44").Value2 = "  -0.07%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.42"
$ws.Range("E45").Value2 = "  -4.07%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.20"
$ws.Range("E46").Value2 = "  -5.76%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.38"
$ws.Range("E47").Value2 = "  -3.67%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.63"
$ws.Range("E48").Value2 = "  -4.52%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.88"
$ws.Range("E49").Value2 = "  -1.01%  "

# Row 50
$ws.Range("D50").Value2 = "2.424.78"
$ws.Range("E50").Value2 = "  -0.65%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.895"
$ws.Range("E51").Value2 = "  -1.98%  "
